$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Update the browser textbox (shape id 27, "CaixaDeTexto 26"): reposition,
#     narrow its height (the autofit engine recomputes Height once the new,
#     single-line text is applied) and replace the two-paragraph caption with
#     a single line of text.
$shp = $s.Shapes.Item(17)
$shp.Left = 10219335 / 12700
$shp.Top = 5911358 / 12700
$shp.Width = 128.08260352519687
$shp.TextFrame.TextRange.Text = "Navegador do usuário"

# --- Add the new "Parte do usuário" label textbox.
$tb1 = $s.Shapes.AddTextbox(1, 8081925 / 12700, 3054046 / 12700, 2049803 / 12700, 338554 / 12700)
$tb1.Name = "CaixaDeTexto 1"
$tb1.Fill.Visible = $false
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1
$tb1.TextFrame.TextRange.Text = "Parte do usuário"
$tb1.TextFrame.TextRange.LanguageID = "pt-BR"
$tb1.TextFrame.TextRange.Font.Size = 16
$tb1.TextFrame.TextRange.Font.Bold = $true
$tb1.TextFrame.TextRange.Font.Name = "Arial"
$tb1.TextFrame2.TextRange.Font.NameComplexScript = "Arial"

# --- Add the new "Parte da solução" label textbox.
$tb2 = $s.Shapes.AddTextbox(1, 2026314 / 12700, 3072357 / 12700, 2097271 / 12700, 338554 / 12700)
$tb2.Name = "CaixaDeTexto 3"
$tb2.Fill.Visible = $false
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1
$tb2.TextFrame.TextRange.Text = "Parte da solução"
$tb2.TextFrame.TextRange.Font.Size = 16
$tb2.TextFrame.TextRange.Font.Bold = $true
$tb2.TextFrame.TextRange.Font.Name = "Arial"
$tb2.TextFrame2.TextRange.Font.NameComplexScript = "Arial"
